$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format for the Price/Volume columns so values such as
# "0.9993" or "1.000" are preserved as literal text (matching the
# original inline-string cell contents) instead of being coerced to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.214.82"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "1.829.17"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "236.73"
$ws.Range("E5").Value = "  -1.36%  "
$ws.Range("D6").Value = "0.6080"
$ws.Range("E6").Value = "  -3.59%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").Value = "0.07102"
$ws.Range("E8").Value = "  -4.98%  "
$ws.Range("D9").Value = "0.2816"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").Value = "23.88"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "0.07669"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").Value = "1.828.13"
$ws.Range("E12").Value = "  -0.77%  "
$ws.Range("D13").Value = "4.829"
$ws.Range("E13").Value = "  -3.04%  "
$ws.Range("D14").Value = "0.00001013"
$ws.Range("E14").Value = "  -1.97%  "
$ws.Range("D15").Value = "0.6377"
$ws.Range("E15").Value = "  -5.92%  "
$ws.Range("D16").Value = "2.074.23"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "79.46"
$ws.Range("E17").Value = "  -3.00%  "
$ws.Range("D18").Value = "5.908"
$ws.Range("E18").Value = "  -5.18%  "
$ws.Range("D19").Value = "29.173.31"
$ws.Range("E19").Value = "  -0.57%  "
$ws.Range("D20").Value = "228.43"
$ws.Range("E20").Value = "  -0.34%  "
$ws.Range("E21").Value = "  -4.18%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "7.033"
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "154.31"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("D26").Value = "8.087"
$ws.Range("E26").Value = "  -5.28%  "
$ws.Range("E27").Value = "  -4.46%  "
$ws.Range("E28").Value = "  -4.71%  "
$ws.Range("D29").Value = "1.495"
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").Value = "0.06483"
$ws.Range("E30").Value = "  -6.07%  "
$ws.Range("E31").Value = "  -2.09%  "
$ws.Range("D32").Value = "3.833"
$ws.Range("E32").Value = "  -5.67%  "
$ws.Range("E33").Value = "  -6.33%  "
$ws.Range("D34").Value = "1.130"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "1.753"
$ws.Range("E35").Value = "  -4.37%  "
$ws.Range("D36").Value = "0.6515"
$ws.Range("E36").Value = "  -7.01%  "
$ws.Range("D37").Value = "2.556"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("D38").Value = "2.757"
$ws.Range("E38").Value = "  -2.16%  "
$ws.Range("D39").Value = "1.218.02"
$ws.Range("E39").Value = "  -1.62%  "
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("D41").Value = "6.512"
$ws.Range("E41").Value = "  -4.33%  "
$ws.Range("D42").Value = "0.9298"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("D44").Value = "101.07"
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "1.979.68"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("E46").Value = "  -3.28%  "
$ws.Range("D47").Value = "0.00000000116"
$ws.Range("E47").Value = "  -2.66%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.585"
$ws.Range("E48").Value = "  -4.38%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.608"
$ws.Range("E49").Value = "  -6.19%  "
$ws.Range("E50").Value = "  -5.93%  "
$ws.Range("D51").Value = "6.488"
